$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 price (SAN939653 / 733739015228)
$ws.Range("C3").Value = 4854.94

# New data rows 4-9: SKU code (A), numeric-looking code (B), price (C)
$newRows = @(
    @("101454368_4725004", "898220010325", 14130.82),
    @("SAN84010",          "306969054093", 6005.6),
    @("SAN0939547",        "728177004613", 9839.4),
    @("100918697_4725004", "898220010332", 6863.54),
    @("SAN83999",          "733739003706", 3162.85),
    @("SAN10",             "733739016539", 14054.09)
)

$row = 4
foreach ($r in $newRows) {
    $ws.Range("A$row").Value = $r[0]

    # Column B values are purely numeric-looking (e.g. "898220010325"); a
    # plain .Value assignment would be auto-typed as a number. Route it
    # through a text-producing formula, then flatten the formula to a
    # static value via copy/paste-values, so it ends up stored as text
    # (matches the workbook's existing convention for these codes) without
    # leaving a formula behind or touching NumberFormat/cell style.
    $bcell = $ws.Range("B$row")
    $bcell.Formula = "=""" + $r[1] + """"
    $bcell.Copy()
    $bcell.PasteSpecial(-4163)  # xlPasteValues

    $ws.Range("C$row").Value = $r[2]
    $row++
}

# Apply the same formatting (style) as the existing data row (row 3) to the
# new rows.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C9").PasteSpecial(-4122)  # xlPasteFormats
